# This edit inserts a new weekly price record for "Betarraga" (row 1275)
# into the existing daily/weekly price log, shifting all subsequent
# records (old rows 1275-1359) down by one row to rows 1276-1360.
#
# The new record:
#   Fecha (D) = 45021  (2023-04-05)
#   Calidad (I) = "Segunda"
#   Volumen (J) = 24000
#   Precio minimo (K) = 60
#   Precio maximo (L) = 60
#   Precio promedio ponderado (M) = 60
#   Precio $/Kg (P) = 60
# All other columns (A,B,C,E,F,G,H,N,O,Q,R) are constant for this
# market/category/variety combination throughout the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 1275; this pushes the former rows
# 1275..1359 down to 1276..1360 (and the sheet dimension grows to R1360).
$ws.Rows.Item(1275).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(1275, 1).Value  = 6
$ws.Cells.Item(1275, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1275, 3).Value  = "Metropolitana"
$ws.Cells.Item(1275, 4).Value  = 45021
$ws.Cells.Item(1275, 5).Value  = 13
$ws.Cells.Item(1275, 6).Value  = 100114014
$ws.Cells.Item(1275, 7).Value  = "Betarraga"
$ws.Cells.Item(1275, 8).Value  = "Sin especificar"
$ws.Cells.Item(1275, 9).Value  = "Segunda"
$ws.Cells.Item(1275, 10).Value = 24000
$ws.Cells.Item(1275, 11).Value = 60
$ws.Cells.Item(1275, 12).Value = 60
$ws.Cells.Item(1275, 13).Value = 60
$ws.Cells.Item(1275, 14).Value = "$/unidad"
$ws.Cells.Item(1275, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1275, 16).Value = 60
$ws.Cells.Item(1275, 17).Value = 1
$ws.Cells.Item(1275, 18).Value = "Hortaliza"
